$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.104.26'
$ws.Range('E2').Value = '  -0.10%  '

$ws.Range('D3').Value = '2.371.16'
$ws.Range('E3').Value = '  +1.32%  '

$ws.Range('E4').Value = '  +0.01%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '303.60'

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '96.28'
$ws.Range('E6').Value = '  +1.26%  '

$ws.Range('E7').Value = '  -0.15%  '

$ws.Range('E8').Value = '  -0.02%  '

$ws.Range('E9').Value = '  -2.61%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '34.42'

$ws.Range('E11').Value = '  +3.95%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0789'
$ws.Range('E12').Value = '  +0.48%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '18.33'
$ws.Range('E13').Value = '  -1.94%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.80'
$ws.Range('E14').Value = '  +0.54%  '

$ws.Range('D15').Value = '2.736.91'
$ws.Range('E15').Value = '  +1.11%  '

$ws.Range('D16').Value = '2.427.07'
$ws.Range('E16').Value = '  +4.01%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.803'
$ws.Range('E17').Value = '  +0.65%  '

$ws.Range('D18').Value = '43.122.12'
$ws.Range('E18').Value = '  +0.09%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.00'
$ws.Range('E19').Value = '  -1.50%  '

$ws.Range('E20').Value = '  +1.55%  '

$ws.Range('D21').Value = '0.0₃0890'
$ws.Range('E21').Value = '  -0.14%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '68.04'
$ws.Range('E22').Value = '  +0.05%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '235.62'
$ws.Range('E23').Value = '  -0.16%  '

$ws.Range('E24').Value = '  +0.38%  '

$ws.Range('E25').Value = '  +1.39%  '

$ws.Range('E26').Value = '  -0.14%  '

$ws.Range('E27').Value = '  -0.65%  '

$ws.Range('E28').Value = '  +4.00%  '

$ws.Range('E29').Value = '  +1.90%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '31.88'
$ws.Range('E30').Value = '  +1.61%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.00'
$ws.Range('E31').Value = '  -0.09%  '

$ws.Range('E32').Value = '  +0.63%  '

$ws.Range('B33').Value = 'Kaspa'
$ws.Range('C33').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.112'
$ws.Range('E33').Value = '  +11.10%  '

$ws.Range('B34').Value = 'Celestia'
$ws.Range('C34').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '17.86'
$ws.Range('E34').Value = '  +3.21%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.0738'
$ws.Range('E35').Value = '  +0.86%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '128.10'
$ws.Range('E36').Value = '  +16.06%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.82'
$ws.Range('E37').Value = '  -0.15%  '

$ws.Range('E38').Value = '  +4.09%  '

$ws.Range('E39').Value = '  -0.56%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.25'
$ws.Range('E40').Value = '  -2.96%  '

$ws.Range('E41').Value = '  -0.70%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '21.22'
$ws.Range('E42').Value = '  -4.42%  '

$ws.Range('D43').Value = '1.929.20'

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0278'
$ws.Range('E44').Value = '  -1.21%  '

$ws.Range('E45').Value = '  +1.74%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.76'
$ws.Range('E46').Value = '  +1.43%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '9.24'
$ws.Range('E47').Value = '  -7.87%  '

$ws.Range('D48').Value = '2.596.52'
$ws.Range('E48').Value = '  +1.01%  '

$ws.Range('E49').Value = '  +3.21%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '71.73'
$ws.Range('E50').Value = '  -0.50%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '51.79'
$ws.Range('E51').Value = '  -2.49%  '
